# Update "Översikt EKERÖ" worksheet:
#  - Row 2 becomes the "A 14020-2023" record (previously row 3), with refreshed
#    survey numbers and two additional species found, and an updated
#    "Förändrad" (changed) date.
#  - Row 3 becomes the "A 8216-2023" record (previously row 2), unchanged
#    apart from the updated "Förändrad" date.
#  - Every other data row (4-45) only gets its "Förändrad" date bumped from
#    2023-09-13 (45182) to 2023-09-15 (45184).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 : A 14020-2023 -------------------------------------------------
$ws.Range("A2").Value = "A 14020-2023"
$ws.Range("B2").Value = 45008
$ws.Range("C2").Value = 45184
$ws.Range("D2").Value = "STOCKHOLMS LÄN"
$ws.Range("E2").Value = "EKERÖ"
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = 9.300000000000001
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 4
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 7
$ws.Range("R2").Value = "Ryl`r`nBacktimjan`r`nOrange taggsvamp`r`nSpillkråka`r`nDropptaggsvamp`r`nGrönpyrola`r`nKopparödla"
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_EKERO/artfynd/A 14020-2023.xlsx")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_EKERO/kartor/A 14020-2023.png")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_EKERO/klagomål/A 14020-2023.docx")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_EKERO/klagomålsmail/A 14020-2023.docx")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_EKERO/tillsyn/A 14020-2023.docx")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_EKERO/tillsynsmail/A 14020-2023.docx")'

# ---- Row 3 : A 8216-2023 ---------------------------------------------------
$ws.Range("A3").Value = "A 8216-2023"
$ws.Range("B3").Value = 44974
$ws.Range("C3").Value = 45184
$ws.Range("D3").Value = "STOCKHOLMS LÄN"
$ws.Range("E3").Value = "EKERÖ"
$ws.Range("F3").Value = "Kommuner"
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 5
$ws.Range("R3").Value = "Gul lammticka`r`nBlomkålssvamp`r`nGranbarkgnagare`r`nThomsons trägnagare`r`nBlåsippa"
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_EKERO/artfynd/A 8216-2023.xlsx")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_EKERO/kartor/A 8216-2023.png")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_EKERO/klagomål/A 8216-2023.docx")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_EKERO/klagomålsmail/A 8216-2023.docx")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_EKERO/tillsyn/A 8216-2023.docx")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/Logging_EKERO/tillsynsmail/A 8216-2023.docx")'

# Re-assigning the wrapped "Artnamn" text above can trigger an autofit of the
# row height; the source file keeps these rows at their original fixed
# height, so restore it explicitly.
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15

# ---- Rows 4-45 : bump "Förändrad" date only -------------------------------
for ($r = 4; $r -le 45; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45182) {
        $cell.Value = 45184
    }
}
